$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.903.42'
$ws.Range("E2").Value = '  -0.19%  '

# Row 3
$ws.Range("D3").Value = '1.875.19'
$ws.Range("E3").Value = '  -0.99%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7373'
$ws.Range("E5").Value = '  -5.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.25'
$ws.Range("E6").Value = '  -0.71%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("E8").Value = '  +0.86%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07189'
$ws.Range("E9").Value = '  -0.84%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.75'
$ws.Range("E10").Value = '  -4.33%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08371'
$ws.Range("E11").Value = '  -3.68%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.906.82'
$ws.Range("E12").Value = '  -9.62%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7503'
$ws.Range("E13").Value = '  -3.20%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.418'
$ws.Range("E14").Value = '  -0.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.61'
$ws.Range("E15").Value = '  -2.09%  '

# Row 16
$ws.Range("D16").Value = '29.915.46'
$ws.Range("E16").Value = '  -1.21%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.073'
$ws.Range("E17").Value = '  -1.80%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '246.50'
$ws.Range("E18").Value = '  +0.09%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.56'
$ws.Range("E19").Value = '  -2.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007836'
$ws.Range("E20").Value = '  -0.37%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9987'
$ws.Range("E21").Value = '  -0.16%  '

# Row 22
$ws.Range("D22").Value = '2.125.97'
$ws.Range("E22").Value = '  -11.22%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.975'
$ws.Range("E23").Value = '  -1.98%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9986'
$ws.Range("E24").Value = '  -0.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1558'
$ws.Range("E25").Value = '  -6.55%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.269'
$ws.Range("E26").Value = '  -2.38%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.96'
$ws.Range("E27").Value = '  +0.96%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.64'
$ws.Range("E28").Value = '  -1.04%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.032'
$ws.Range("E29").Value = '  -1.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.510'
$ws.Range("E30").Value = '  +5.30%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.594'
$ws.Range("E31").Value = '  +1.58%  '

# Row 32
$ws.Range("E32").Value = '  -0.72%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.269'
$ws.Range("E33").Value = '  +3.35%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05321'
$ws.Range("E34").Value = '  -2.79%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.239'
$ws.Range("E35").Value = '  -0.54%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7523'
$ws.Range("E36").Value = '  -0.56%  '

# Row 37
$ws.Range("E37").Value = '  -0.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.697'
$ws.Range("E38").Value = '  +0.06%  '

# Row 39
$ws.Range("E39").Value = '  -0.71%  '

# Row 40
$ws.Range("E40").Value = '  -1.26%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4504'
$ws.Range("E41").Value = '  -0.17%  '

# Row 42
$ws.Range("D42").Value = '1.112.65'
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.054'
$ws.Range("E43").Value = '  -1.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.54'
$ws.Range("E44").Value = '  -1.53%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8544'
$ws.Range("E45").Value = '  +0.52%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.01%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.10'
$ws.Range("E47").Value = '  -0.68%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.639'
$ws.Range("E48").Value = '  +0.28%  '

# Row 49
$ws.Range("B49").Value = 'SynthetixNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.100'
$ws.Range("E49").Value = '  +2.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.842'
$ws.Range("E50").Value = '  -1.95%  '

# Row 51
$ws.Range("D51").Value = '2.023.12'
$ws.Range("E51").Value = '  -9.16%  '
